$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp in the title cell (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 15:34"

# --- Row 4: Estados Unidos (covid figures refresh) ---
$ws.Range("B4").Value = 1323077
$ws.Range("C4").Value = 1292
$ws.Range("D4").Value = 223876
$ws.Range("E4").Value = 1020564
$ws.Range("F4").Value = 16917
$ws.Range("G4").Value = 22
$ws.Range("H4").Value = 78637

# --- Row 20: Arabia Saudita (covid figures refresh) ---
$ws.Range("B20").Value = 37136
$ws.Range("C20").Value = 1704
$ws.Range("D20").Value = 10144
$ws.Range("E20").Value = 26753
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 239

# --- Row 45: Serbia (covid figures refresh) ---
$ws.Range("B45").Value = 10032
$ws.Range("C45").Value = 89
$ws.Range("D45").Value = 2732
$ws.Range("E45").Value = 7087
$ws.Range("F45").Value = 43
$ws.Range("G45").Value = 4
$ws.Range("H45").Value = 213

# --- Row 103: Sri Lanka (covid figures refresh) ---
$ws.Range("D103").Value = 260
$ws.Range("E103").Value = 575

# --- Rows 135-139: Benin's total rises above Congo/Ruanda/Chad/Zambia, so it
#     re-sorts to the top of this block (right after Vietnam); the other four
#     countries shift down one row each, keeping their own figures unchanged. ---
$ws.Range("A135").Value = "Benin"
$ws.Range("B135").Value = 284
$ws.Range("C135").Value = 42
$ws.Range("D135").Value = 62
$ws.Range("E135").Value = 220
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 2

$ws.Range("A136").Value = "Congo"
$ws.Range("B136").Value = 274
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 33
$ws.Range("E136").Value = 231
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 10

$ws.Range("A137").Value = "Ruanda"
$ws.Range("B137").Value = 273
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 136
$ws.Range("E137").Value = 137
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 0

$ws.Range("A138").Value = "Republica del Chad"
$ws.Range("B138").Value = 260
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 50
$ws.Range("E138").Value = 182
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 28

$ws.Range("A139").Value = "Zambia"
$ws.Range("B139").Value = 252
$ws.Range("C139").Value = 85
$ws.Range("D139").Value = 112
$ws.Range("E139").Value = 133
$ws.Range("F139").Value = 1
$ws.Range("G139").Value = 3
$ws.Range("H139").Value = 7
